# "Error Calculations and Plots" -- refresh this sample's missing-data mask:
# two rows that used to carry a manually-injected "missing" marker (RM 232's
# column D value and SC 92 in full) are dropped from the sheet entirely, and
# the set of cells masked out in column D (header "D" / worksheet column E)
# is rebalanced: a few cells that were blank get a restored numeric reading,
# while a couple that had a reading get re-masked as blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "RM 232" row and the "SC 92" row (identified by their ID in
# column A at the time of edit: rows 26 and 28). Remove the higher-numbered
# row first so the lower row index still points at the right row afterwards.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Rebalance which column-D ("E" in the sheet) readings are masked as missing.
$ws.Range("E3").Value = -5.7
$ws.Range("E5").Value = ""
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").Value = ""
$ws.Range("E32").Value = -6.4
